# Fruta / hortaliza, semanal
# Weekly data refresh: two new price observations were added to the
# "Chirimoya" dataset. Inserting the rows shifts the remaining (unchanged)
# rows down automatically, matching the target layout (dimension A1:T32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row inserted at row 4 -------------------------------------------
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = "7/27/2023"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("Q4").Value = "`$/bandeja 8 kilos"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 3500
$ws.Range("T4").Value = 8

# --- New row inserted at row 23 (post first insert) ----------------------
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = "12/15/2023"
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = "Chirimoya"
$ws.Range("K23").Value = "Cultivar IV Región"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 24000
$ws.Range("P23").Value = 24000
$ws.Range("Q23").Value = "`$/bandeja 8 kilos"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 3000
$ws.Range("T23").Value = 8
